$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new status row (row 20) for the new update entry, reusing the
# date formatting of the row above it (style index 1 - mm/dd/yyyy)
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A20").Value = 43168

$ws.Range("D20").Value = "WebApp edit entry (70%)"
$ws.Range("B20").Value = "Android bugfix post methods -> research, because they don't work (20%)"
$ws.Range("C20").Value = "C# Admin note (100%). Login status bar does not work as it should work (20%)"

$ws.Range("D24").Select()
